$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# The URL that A2's hyperlink points to (and keeps pointing to).
$serverUrl = "http://172.191.4.85/TestCollection"

# Remove every hyperlink on the sheet - in this host, Hyperlinks.Delete()
# clears the whole sheet collection rather than just the scoped range, so
# this also drops A3's hyperlink, which is what we want since row 3 is
# being cleared entirely.
$ws.Range("A1:C3").Hyperlinks.Delete()

# Re-create the hyperlink that should remain on A2.
$ws.Hyperlinks.Add($ws.Range("A2"), $serverUrl)

# Clear the now-unused values in B2 and C2 (keep A2's URL + hyperlink + style)
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()

# Clear row 3 entirely (values only, keep formatting/styles)
$ws.Range("A3").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("C3").ClearContents()

# Update the active selection to C2
$ws.Range("C2").Select()

$wb.Save()
